# Mise à jour de l'application
# Adds 5 new "Entrainement" rows (J-2, Global, 2025-09-04) to the tracking
# sheet, mirroring the layout of the existing data rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 429
$newRowsCount = 5

# Copy the number/alignment formatting from the last existing data row so the
# new rows keep the same look (date format in column B, centered style in
# column D) without introducing any new cell styles.
$srcRange = $ws.Range("A" + $lastRow + ":V" + $lastRow)
$dstRange = $ws.Range("A" + ($lastRow + 1) + ":V" + ($lastRow + $newRowsCount))
$srcRange.Copy()
$dstRange.PasteSpecial(-4122) # xlPasteFormats

$rows = @(
    @{ A = "Entrainement"; B = 45904; C = "Global"; D = "J-2"; E = "Omar Benyounes";     F = "center midfield"; G = "01:29:14";
       H = 5.3;  I = 0.44; J = 4.85; K = 0.27; L = 0.14; M = 0.03; N = 0; O = 5; P = 3.49; Q = 27.12; R = 4.17; S = 28; T = 3; U = 9;  V = 2 },
    @{ A = "Entrainement"; B = 45904; C = "Global"; D = "J-2"; E = "Karahali Souaré";    F = "right forward";   G = "01:26:29";
       H = 4.25; I = 0.18; J = 4.06; K = 0.17; L = 0.01; M = 0;    N = 0; O = 0; P = 2.63; Q = 21.28; R = 5.25; S = 28; T = 9; U = 29; V = 7 },
    @{ A = "Entrainement"; B = 45904; C = "Global"; D = "J-2"; E = "Ilan Ihaddadene";    F = "center midfield"; G = "00:57:08";
       H = 3.23; I = 0.12; J = 3.11; K = 0.12; L = 0.01; M = 0;    N = 0; O = 0; P = 3.26; Q = 22.03; R = 4.06; S = 13; T = 1; U = 7;  V = 0 },
    @{ A = "Entrainement"; B = 45904; C = "Global"; D = "J-2"; E = "Mattheo Haon";       F = "right back";      G = "00:58:30";
       H = 3.45; I = 0.28; J = 3.17; K = 0.09; L = 0.14; M = 0.04; N = 0; O = 5; P = 3.43; Q = 27.43; R = 4.37; S = 8;  T = 4; U = 4;  V = 0 },
    @{ A = "Entrainement"; B = 45904; C = "Global"; D = "J-2"; E = "Ilyes Boughanmi";    F = "center forward";  G = "01:27:53";
       H = 4.18; I = 0.14; J = 4.03; K = 0.14; L = 0;    M = 0;    N = 0; O = 0; P = 2.76; Q = 19.64; R = 4.77; S = 16; T = 8; U = 4;  V = 1 }
)

$r = $lastRow
foreach ($row in $rows) {
    $r = $r + 1

    $ws.Cells.Item($r, 1).Value = $row.A
    $ws.Cells.Item($r, 2).Value = $row.B
    $ws.Cells.Item($r, 3).Value = $row.C
    $ws.Cells.Item($r, 4).Value = $row.D
    $ws.Cells.Item($r, 5).Value = $row.E
    $ws.Cells.Item($r, 6).Value = $row.F
    $ws.Cells.Item($r, 7).Value = $row.G
    $ws.Cells.Item($r, 8).Value = $row.H
    $ws.Cells.Item($r, 9).Value = $row.I
    $ws.Cells.Item($r, 10).Value = $row.J
    $ws.Cells.Item($r, 11).Value = $row.K
    $ws.Cells.Item($r, 12).Value = $row.L
    $ws.Cells.Item($r, 13).Value = $row.M
    $ws.Cells.Item($r, 14).Value = $row.N
    $ws.Cells.Item($r, 15).Value = $row.O
    $ws.Cells.Item($r, 16).Value = $row.P
    $ws.Cells.Item($r, 17).Value = $row.Q
    $ws.Cells.Item($r, 18).Value = $row.R
    $ws.Cells.Item($r, 19).Value = $row.S
    $ws.Cells.Item($r, 20).Value = $row.T
    $ws.Cells.Item($r, 21).Value = $row.U
    $ws.Cells.Item($r, 22).Value = $row.V
}

# Update the saved view: scroll position and the active selection, as left
# by the author after entering the new rows.
$win = $excel.ActiveWindow
$win.ScrollRow = 403
$win.ScrollColumn = 1
$ws.Range("C437").Select() | Out-Null
